# ajout de mon nom
#
# The paragraph "Louis Pelletier & Charlotte _" gets its trailing
# placeholder underscore replaced with "de Lanauze", so the run is split
# into:
#   "Louis Pelletier & Charlotte "   (existing run, now ending in a space)
#   "de Lanauze"                      (new run, identical character formatting)

$d = $word.ActiveDocument

# Find the exact placeholder text and capture the matched Range.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Louis Pelletier & Charlotte _", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    Write-Output "ERROR: placeholder paragraph not found"
} else {
    # The trailing "_" is the final character of the matched range.
    $underscoreStart = $searchRange.End - 1
    $underscoreEnd = $searchRange.End
    $r = $d.Range($underscoreStart, $underscoreEnd)

    # Replace just the "_" with a brand-new run "de Lanauze" that carries
    # the same explicit run formatting (bold, Times New Roman, 32 half-pt,
    # fr-CA) as the run it is being inserted into, so it renders identically
    # to the surrounding text while remaining its own <w:r>.
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="fr-CA"/></w:rPr><w:t>de Lanauze</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)

    Write-Output "Updated paragraph text: $($d.Paragraphs.Item(9).Range.Text)"
}
